$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fixture")

# Row 1: Manchester City v Leicester City -> Swindon Town v Manchester City
#        26/12/2021 15:00 | Premier League -> 07/01/2022 20:00 | FA CUP
$ws.Range("A1").Value = "Swindon Town v Manchester City"
$ws.Range("B1").Value = "07/01/2022 20:00 | FA CUP"

# Row 4: kickoff time change only
#        19/02/2022 15:00 | Premier League -> 19/02/2022 17:30 | Premier League
$ws.Range("B4").Value = "19/02/2022 17:30 | Premier League"
